$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.028.79"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.786.58"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.04"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4541"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3602"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07504"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.26"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9970"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.03"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.069"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.239"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "1.781.96"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.86"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001062"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06438"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9986"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.822"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "28.074.56"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.083"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.90"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.43"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "1.984.45"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.250"
$ws.Range("E29").Value = "  +7.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.26"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.116"
$ws.Range("E31").Value = "  +3.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09183"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.676"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.586"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.95"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02304"
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06153"
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2098"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6365"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.990"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.189"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.391"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.921"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.28"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5941"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.740"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.09"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.970"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06952"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.139"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.06"
$ws.Range("E51").Value = "  +0.80%  "
